# Fruta / hortaliza, semanal
# Insert two new weekly records at rows 89-90 (pushing all existing
# records at/after the old row 89 down by two rows), then populate the
# two new rows with the new "Provincia de Limarí" observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 89 (Excel shifts rows
# 89:175 down to 91:177 and carries the formatting of the row above,
# which is what we want for the date-formatted column D).
$ws.Rows.Item(89).Insert()
$ws.Rows.Item(89).Insert()

# New row 89
$ws.Cells.Item(89, 1).Value = 8
$ws.Cells.Item(89, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(89, 3).Value = "Coquimbo"
$ws.Cells.Item(89, 4).Value2 = 45271
$ws.Cells.Item(89, 5).Value = 4
$ws.Cells.Item(89, 6).Value = 100112028
$ws.Cells.Item(89, 7).Value = "Sandia"
$ws.Cells.Item(89, 8).Value = "Sin especificar"
$ws.Cells.Item(89, 9).Value = "Primera"
$ws.Cells.Item(89, 10).Value = 2000
$ws.Cells.Item(89, 11).Value = 650
$ws.Cells.Item(89, 12).Value = 700
$ws.Cells.Item(89, 13).Value = 675
$ws.Cells.Item(89, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(89, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(89, 16).Value = 675
$ws.Cells.Item(89, 17).Value = 1
$ws.Cells.Item(89, 18).Value = "Hortaliza"

# New row 90
$ws.Cells.Item(90, 1).Value = 8
$ws.Cells.Item(90, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(90, 3).Value = "Coquimbo"
$ws.Cells.Item(90, 4).Value2 = 45271
$ws.Cells.Item(90, 5).Value = 4
$ws.Cells.Item(90, 6).Value = 100112028
$ws.Cells.Item(90, 7).Value = "Sandia"
$ws.Cells.Item(90, 8).Value = "Sin especificar"
$ws.Cells.Item(90, 9).Value = "Segunda"
$ws.Cells.Item(90, 10).Value = 1400
$ws.Cells.Item(90, 11).Value = 500
$ws.Cells.Item(90, 12).Value = 550
$ws.Cells.Item(90, 13).Value = 525
$ws.Cells.Item(90, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(90, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(90, 16).Value = 525
$ws.Cells.Item(90, 17).Value = 1
$ws.Cells.Item(90, 18).Value = "Hortaliza"
